# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" quarterly sheet (right after "总计") and updates the
# "总计" (summary) sheet with the new quarter's totals, shifting the older
# quarters down by one row - per commit "feat: add 2022-Q3 data".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 0. Grab a handle on the existing "2022-Q2" sheet *before* we insert the new
#    sheet - we reuse its already-correct cell formatting (bold/border/center
#    header style, and the bold/border/center index-column style) instead of
#    re-deriving the look from scratch, so the new sheet matches its siblings
#    exactly.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$refSheet = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet immediately after "总计".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# Header row (B1:H1).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Fund holdings data, one row per fund, ordered exactly as in the source.
# Columns: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元), 仓位排名
$q3Data = @(
    @("001184","易方达新常态灵活配置混合","22.40","85.86","3.29","0.7370",10),
    @("159745","国泰中证全指建筑材料ETF","7.92","99.14","7.79","0.6170",3),
    @("004856","广发中证全指建筑材料指数A","7.66","93.74","7.36","0.5638",3),
    @("004857","广发中证全指建筑材料指数C","6.12","93.74","7.36","0.4504",3),
    @("012437","德邦价值优选混合A","6.73","90.21","3.46","0.2329",10),
    @("420005","天弘周期策略混合A","3.10","93.08","6.18","0.1916",6),
    @("001179","德邦大健康灵活配置混合","3.96","89.58","4.00","0.1584",6),
    @("516750","富国中证全指建筑材料ETF","0.82","98.46","7.84","0.0643",3),
    @("012419","天弘国证建材指数C","0.63","94.93","7.26","0.0457",4),
    @("015458","天弘周期策略混合C","0.68","93.08","6.18","0.0420",6),
    @("008840","德邦大消费混合A","1.06","90.05","3.89","0.0412",7),
    @("006167","德邦乐享生活混合A","1.03","90.49","2.87","0.0296",9),
    @("008841","德邦大消费混合C","0.56","90.05","3.89","0.0218",7),
    @("159787","易方达中证全指建筑材料ETF","0.17","94.24","7.48","0.0127",3),
    @("006168","德邦乐享生活混合C","0.40","90.49","2.87","0.0115",9),
    @("012438","德邦价值优选混合C","0.28","90.21","3.46","0.0097",10),
    @("012405","天弘国证建材指数A","0.13","94.93","7.26","0.0094",4)
)

# Fund code (B) and the D:G metrics are stored as literal text in the source
# workbook (e.g. "22.40", not 22.4), so force text format on those columns
# before assigning - otherwise Excel auto-coerces them to numbers and mangles
# leading zeros / trailing zeros.
$lastRow = 1 + $q3Data.Count
$q3.Range("B2:B" + $lastRow).NumberFormat = "@"
$q3.Range("D2:G" + $lastRow).NumberFormat = "@"

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    $q3.Cells.Item($r, 2).Value = $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = $row[2]
    $q3.Cells.Item($r, 5).Value = $row[3]
    $q3.Cells.Item($r, 6).Value = $row[4]
    $q3.Cells.Item($r, 7).Value = $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Reuse the header style (bold, centered, bordered) and the index-column
# style from the sibling "2022-Q2" sheet so the new sheet is visually
# identical to the others.
$refSheet.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$refSheet.Range("A2:A" + $lastRow).Copy()
$q3.Range("A2:A" + $lastRow).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: a new top data row for 2022-Q3, and every older
#    quarter shifts down one row (with its running index bumped by one).
# ---------------------------------------------------------------------------
$totals = @(
    @("2022-Q3", 17, 3.24),
    @("2022-Q2", 40, 5.32),
    @("2022-Q1", 21, 6.82),
    @("2021-Q4", 54, 12.55),
    @("2021-Q3", 30, 23.14),
    @("2021-Q2", 32, 29.74),
    @("2021-Q1", 45, 39.05),
    @("2020-Q4", 37, 36.01)
)

$r = 2
foreach ($row in $totals) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
    $totalSheet.Cells.Item($r, 2).Value = $row[0]
    $totalSheet.Cells.Item($r, 3).Value = $row[1]
    $totalSheet.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}
